$d = $word.ActiveDocument

$tags = @(
  "PUMP:HRD:0000",
  "PUMP:HRS:3340",
  "PUMP:HTR:100",
  "PUMP:HTR:1100",
  "PUMP:HTR:1200",
  "PUMP:HTR:1300",
  "PUMP:HTR:1400",
  "PUMP:HTR:1500",
  "PUMP:HTR:300",
  "PUMP:HTR:400",
  "PUMP:HTR:500",
  "PUMP:INS:100",
  "PUMP:INS:110",
  "PUMP:INS:120",
  "PUMP:INS:130",
  "PUMP:INS:140",
  "PUMP:INS:150",
  "PUMP:INS:160",
  "PUMP:INS:170",
  "PUMP:INS:180",
  "PUMP:INS:190",
  "PUMP:INS:200",
  "PUMP:INS:210",
  "PUMP:INS:220",
  "PUMP:PRS:2",
  "PUMP:PRS:3340",
  "PUMP:SDS:10",
  "PUMP:SDS:20",
  "PUMP:SDS:30",
  "PUMP:SDS:40",
  "PUMP:SDS:50",
  "PUMP:SDS:60",
  "PUMP:SDS:70",
  "PUMP:SVATR:100",
  "PUMP:SVATR:200",
  "PUMP:SVATR:300",
  "PUMP:SVATR:400",
  "PUMP:SVATR:500",
  "PUMP:UT:100",
  "PUMP:UT:110",
  "PUMP:UT:120",
  "PUMP:UT:130",
  "PUMP:UT:140",
  "PUMP:UT:150",
  "PUMP:UT:160",
  "PUMP:UT:170",
  "PUMP:UT:180",
  "PUMP:UT:190",
  "PUMP:UT:200",
  "PUMP:UT:210",
  "PUMP:UT:220",
  "PUMPHTR:200"
)

# Use the (non-bold) Title paragraph's own paragraph mark as a clean,
# unformatted template. Copying just that one mark character onto the
# very end of the document splits off a brand-new, completely bare
# paragraph (no rPr/pPr baggage at all) each time - unlike
# InsertParagraphAfter(), which clones the bold run that precedes it.
$titlePara = $d.Paragraphs(1).Range
$cleanMark = $d.Range($titlePara.End - 1, $titlePara.End)

foreach ($tag in $tags) {
  $endPos = $d.Content.End
  $insertion = $d.Range($endPos, $endPos)
  $insertion.FormattedText = $cleanMark.FormattedText

  $newPara = $d.Paragraphs.Last
  $newPara.Range.Text = $tag
}

Write-Output "done"
